# Project Summary.docx update
# - "user friendly" -> "user-friendly" (plus a collapsed _GoBack bookmark
#   at the point where the replacement text was typed, as Word leaves behind)
# - removed a stray leading double-space before "Options to share..."
# - fixed "Te idea" -> "The idea" typo
# - fixed "together( we" -> "together (we" (misplaced space before the
#   parenthesis)
# - fixed "likeh k-nearest" -> "like k-nearest" (stray extra "h")
# - inserted three additional blank paragraphs right before "Future Works:"

$d = $word.ActiveDocument

# --- "Straightforward and user friendly GUI..." -> hyphenate "user-friendly"
# Split the run right after "Straightforward and " (this is also where Word
# drops its "_GoBack" last-edit bookmark) and then fix the wording.
$rng = $d.Content
$found = $rng.Find.Execute("Straightforward and ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $point = $rng.Duplicate
    $point.Start = $rng.End
    $point.End = $rng.End
    $d.Bookmarks.Add("_GoBack", $point)
}
$d.Content.Find.Execute("user friendly", $true, $false, $false, $false, $false, $true, 1, $false, "user-friendly", 2)

# --- remove the stray leading double space before "Options to share..."
$d.Content.Find.Execute("  Options to share", $true, $false, $false, $false, $false, $true, 1, $false, "Options to share", 2)

# --- fix typo: "Te idea" -> "The idea"
$d.Content.Find.Execute("Te idea behind CLIP", $true, $false, $false, $false, $false, $true, 1, $false, "The idea behind CLIP", 2)

# --- fix misplaced space: "together( we" -> "together (we"
$d.Content.Find.Execute("together( we", $true, $false, $false, $false, $false, $true, 1, $false, "together (we", 2)

# --- fix typo: "likeh k-nearest" -> "like k-nearest"
$d.Content.Find.Execute("likeh k-nearest", $true, $false, $false, $false, $false, $true, 1, $false, "like k-nearest", 2)

# --- insert three blank paragraphs between "All we need..." and "Future Works:"
$tail = $d.Content
$tailFound = $tail.Find.Execute("All we need for the completion", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($tailFound) {
    $equipPara = $tail.Paragraphs(1)
    $equipPara.Range.InsertParagraphAfter()
    $equipPara.Range.InsertParagraphAfter()
    $equipPara.Range.InsertParagraphAfter()
}
